$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 1444.3334  # H70: 1396.6666 -> 1444.3334
$ws.Cells.Item(70, 9).Value = 1000  # I70: 1066.6666 -> 1000
$ws.Cells.Item(70, 10).Value = 1571.2858  # J70: 1506.6666 -> 1571.2858
$ws.Cells.Item(70, 11).Value = 3000  # K70: 3199.9998 -> 3000
$ws.Cells.Item(70, 12).Value = 4713.857400000001  # L70: 4519.9998 -> 4713.857400000001
$ws.Cells.Item(70, 13).Value = -2730  # M70: -2929.9998 -> -2730
$ws.Cells.Item(70, 14).Value = -5253.857400000001  # N70: -5059.9998 -> -5253.857400000001

$ws.Cells.Item(73, 8).Value = 1444.3334  # H73: 1396.6666 -> 1444.3334
$ws.Cells.Item(73, 9).Value = 1000  # I73: 1066.6666 -> 1000
$ws.Cells.Item(73, 10).Value = 1571.2858  # J73: 1506.6666 -> 1571.2858
$ws.Cells.Item(73, 11).Value = 3000  # K73: 3199.9998 -> 3000
$ws.Cells.Item(73, 12).Value = 4713.857400000001  # L73: 4519.9998 -> 4713.857400000001
$ws.Cells.Item(73, 13).Value = -2064  # M73: -2263.9998 -> -2064
$ws.Cells.Item(73, 14).Value = -6585.857400000001  # N73: -6391.9998 -> -6585.857400000001

$ws.Cells.Item(75, 8).Value = 29999.5  # H75: 33000 -> 29999.5
$ws.Cells.Item(75, 10).Value = 29999.5  # J75: 33000 -> 29999.5
$ws.Cells.Item(75, 12).Value = 29999.5  # L75: 33000 -> 29999.5
$ws.Cells.Item(75, 14).Value = -31871.5  # N75: -34872 -> -31871.5

$ws.Cells.Item(78, 8).Value = 29999.5  # H78: 33000 -> 29999.5
$ws.Cells.Item(78, 10).Value = 29999.5  # J78: 33000 -> 29999.5
$ws.Cells.Item(78, 12).Value = 89998.5  # L78: 99000 -> 89998.5
$ws.Cells.Item(78, 14).Value = -99358.5  # N78: -108360 -> -99358.5

$ws.Cells.Item(100, 8).Value = 1454.6154  # H100: 1573.4546 -> 1454.6154
$ws.Cells.Item(100, 9).Value = 1454.6154  # I100: 1328.8572 -> 1454.6154
$ws.Cells.Item(100, 10).Value = 0  # J100: 2001.5 -> 0
$ws.Cells.Item(100, 11).Value = 1454.6154  # K100: 1328.8572 -> 1454.6154
$ws.Cells.Item(100, 12).Value = 0  # L100: 2001.5 -> 0
$ws.Cells.Item(100, 13).Value = -913.6153999999999  # M100: -787.8571999999999 -> -913.6153999999999
$ws.Cells.Item(100, 14).ClearContents()  # N100: remove (was -3083.5)

$ws.Cells.Item(116, 8).Value = 2501.6667  # H116: 0 -> 2501.6667
$ws.Cells.Item(116, 9).Value = 2501.6667  # I116: 0 -> 2501.6667
$ws.Cells.Item(116, 11).Value = 2501.6667  # K116: 0 -> 2501.6667
$ws.Cells.Item(116, 13).Value = 940.3332999999998  # M116: None -> 940.3332999999998

$ws.Cells.Item(132, 8).Value = 1715.8214  # H132: 2025.1132 -> 1715.8214
$ws.Cells.Item(132, 9).Value = 1501.6111  # I132: 1696.5 -> 1501.6111
$ws.Cells.Item(132, 10).Value = 7499.5  # J132: 5179.8 -> 7499.5
$ws.Cells.Item(132, 11).Value = 4504.8333  # K132: 5089.5 -> 4504.8333
$ws.Cells.Item(132, 12).Value = 22498.5  # L132: 15539.4 -> 22498.5
$ws.Cells.Item(132, 13).Value = -1974.8333  # M132: -2559.5 -> -1974.8333
$ws.Cells.Item(132, 14).Value = -27558.5  # N132: -20599.4 -> -27558.5

$ws.Cells.Item(138, 8).Value = 1815.3962  # H138: 1876 -> 1815.3962
$ws.Cells.Item(138, 10).Value = 3824.2  # J138: 3621.6924 -> 3824.2
$ws.Cells.Item(138, 12).Value = 11472.6  # L138: 10865.0772 -> 11472.6
$ws.Cells.Item(138, 14).Value = -21752.6  # N138: -21145.0772 -> -21752.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(62, 8).Value = 30000  # H62: 0 -> 30000
$ws.Cells.Item(62, 10).Value = 30000  # J62: 0 -> 30000
$ws.Cells.Item(62, 12).Value = 30000  # L62: 0 -> 30000
$ws.Cells.Item(62, 14).Value = -31248  # N62: None -> -31248

$ws.Cells.Item(65, 8).Value = 30000  # H65: 0 -> 30000
$ws.Cells.Item(65, 10).Value = 30000  # J65: 0 -> 30000
$ws.Cells.Item(65, 12).Value = 90000  # L65: 0 -> 90000
$ws.Cells.Item(65, 14).Value = -96240  # N65: None -> -96240

$ws.Cells.Item(70, 8).Value = 100000  # H70: 0 -> 100000
$ws.Cells.Item(70, 10).Value = 100000  # J70: 0 -> 100000
$ws.Cells.Item(70, 12).Value = 100000  # L70: 0 -> 100000
$ws.Cells.Item(70, 14).Value = -100540  # N70: None -> -100540

$ws.Cells.Item(73, 8).Value = 100000  # H73: 0 -> 100000
$ws.Cells.Item(73, 10).Value = 100000  # J73: 0 -> 100000
$ws.Cells.Item(73, 12).Value = 100000  # L73: 0 -> 100000
$ws.Cells.Item(73, 14).Value = -101872  # N73: None -> -101872

$ws.Cells.Item(74, 8).Value = 2450.861  # H74: 2684.6785 -> 2450.861
$ws.Cells.Item(74, 9).Value = 2190.16  # I74: 2334.7 -> 2190.16
$ws.Cells.Item(74, 10).Value = 3043.3635  # J74: 3559.625 -> 3043.3635
$ws.Cells.Item(74, 11).Value = 2190.16  # K74: 2334.7 -> 2190.16
$ws.Cells.Item(74, 12).Value = 3043.3635  # L74: 3559.625 -> 3043.3635
$ws.Cells.Item(74, 13).Value = -1316.16  # M74: -1460.7 -> -1316.16
$ws.Cells.Item(74, 14).Value = -4791.363499999999  # N74: -5307.625 -> -4791.363499999999

$ws.Cells.Item(77, 8).Value = 2450.861  # H77: 2684.6785 -> 2450.861
$ws.Cells.Item(77, 9).Value = 2190.16  # I77: 2334.7 -> 2190.16
$ws.Cells.Item(77, 10).Value = 3043.3635  # J77: 3559.625 -> 3043.3635
$ws.Cells.Item(77, 11).Value = 10950.8  # K77: 11673.5 -> 10950.8
$ws.Cells.Item(77, 12).Value = 15216.8175  # L77: 17798.125 -> 15216.8175
$ws.Cells.Item(77, 13).Value = -6582.799999999999  # M77: -7305.5 -> -6582.799999999999
$ws.Cells.Item(77, 14).Value = -23952.8175  # N77: -26534.125 -> -23952.8175

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(93, 8).Value = 59500  # H93: 61500 -> 59500
$ws.Cells.Item(93, 10).Value = 59500  # J93: 61500 -> 59500
$ws.Cells.Item(93, 12).Value = 59500  # L93: 61500 -> 59500
$ws.Cells.Item(93, 14).Value = -63244  # N93: -65244 -> -63244

$ws.Cells.Item(133, 8).Value = 55690  # H133: 55780 -> 55690
$ws.Cells.Item(133, 10).Value = 55690  # J133: 55780 -> 55690
$ws.Cells.Item(133, 12).Value = 55690  # L133: 55780 -> 55690
$ws.Cells.Item(133, 14).Value = -65810  # N133: -65900 -> -65810

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 6242.9434  # H31: 8458.433000000001 -> 6242.9434
$ws.Cells.Item(31, 9).Value = 1334.25  # I31: 1572 -> 1334.25
$ws.Cells.Item(31, 10).Value = 10305.311  # J31: 12650.174 -> 10305.311
$ws.Cells.Item(31, 11).Value = 1334.25  # K31: 1572 -> 1334.25
$ws.Cells.Item(31, 12).Value = 10305.311  # L31: 12650.174 -> 10305.311
$ws.Cells.Item(31, 13).Value = -1039.25  # M31: -1277 -> -1039.25
$ws.Cells.Item(31, 14).Value = -10895.311  # N31: -13240.174 -> -10895.311

$ws.Cells.Item(34, 8).Value = 6242.9434  # H34: 8458.433000000001 -> 6242.9434
$ws.Cells.Item(34, 9).Value = 1334.25  # I34: 1572 -> 1334.25
$ws.Cells.Item(34, 10).Value = 10305.311  # J34: 12650.174 -> 10305.311
$ws.Cells.Item(34, 11).Value = 1334.25  # K34: 1572 -> 1334.25
$ws.Cells.Item(34, 12).Value = 10305.311  # L34: 12650.174 -> 10305.311
$ws.Cells.Item(34, 13).Value = -1132.25  # M34: -1370 -> -1132.25
$ws.Cells.Item(34, 14).Value = -10709.311  # N34: -13054.174 -> -10709.311

$ws.Cells.Item(110, 8).Value = 41111  # H110: 0 -> 41111
$ws.Cells.Item(110, 10).Value = 41111  # J110: 0 -> 41111
$ws.Cells.Item(110, 12).Value = 41111  # L110: 0 -> 41111
$ws.Cells.Item(110, 14).Value = -49291  # N110: None -> -49291

$ws.Cells.Item(118, 8).Value = 0  # H118: 38333 -> 0
$ws.Cells.Item(118, 10).Value = 0  # J118: 38333 -> 0
$ws.Cells.Item(118, 12).Value = 0  # L118: 38333 -> 0
$ws.Cells.Item(118, 14).ClearContents()  # N118: remove (was -41647)

$ws.Cells.Item(132, 8).Value = 14958165  # H132: 15152429 -> 14958165
$ws.Cells.Item(132, 9).Value = 19608658  # I132: 19608662 -> 19608658
$ws.Cells.Item(132, 10).Value = 6173902.5  # J132: 6411355.5 -> 6173902.5
$ws.Cells.Item(132, 11).Value = 58825974  # K132: 58825986 -> 58825974
$ws.Cells.Item(132, 12).Value = 18521707.5  # L132: 19234066.5 -> 18521707.5
$ws.Cells.Item(132, 13).Value = -58823444  # M132: -58823456 -> -58823444
$ws.Cells.Item(132, 14).Value = -18526767.5  # N132: -19239126.5 -> -18526767.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(131, 8).Value = 6915.6  # H131: 5400.4443 -> 6915.6
$ws.Cells.Item(131, 9).Value = 445  # I131: 435.83334 -> 445
$ws.Cells.Item(131, 10).Value = 10399.77  # J131: 6818.905 -> 10399.77
$ws.Cells.Item(131, 11).Value = 1335  # K131: 1307.50002 -> 1335
$ws.Cells.Item(131, 12).Value = 31199.31  # L131: 20456.715 -> 31199.31
$ws.Cells.Item(131, 13).Value = 3705  # M131: 3732.49998 -> 3705
$ws.Cells.Item(131, 14).Value = -41279.31  # N131: -30536.715 -> -41279.31

$ws.Cells.Item(141, 8).Value = 7484.125  # H141: 7466.0435 -> 7484.125
$ws.Cells.Item(141, 9).Value = 2969.5557  # I141: 2340.75 -> 2969.5557
$ws.Cells.Item(141, 10).Value = 10192.866  # J141: 10199.533 -> 10192.866
$ws.Cells.Item(141, 11).Value = 8908.667099999999  # K141: 7022.25 -> 8908.667099999999
$ws.Cells.Item(141, 12).Value = 30578.598  # L141: 30598.599 -> 30578.598
$ws.Cells.Item(141, 13).Value = -3728.667099999999  # M141: -1842.25 -> -3728.667099999999
$ws.Cells.Item(141, 14).Value = -40938.598  # N141: -40958.599 -> -40938.598

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(94, 8).Value = 95448  # H94: 45198 -> 95448
$ws.Cells.Item(94, 10).Value = 95448  # J94: 45198 -> 95448
$ws.Cells.Item(94, 12).Value = 95448  # L94: 45198 -> 95448
$ws.Cells.Item(94, 14).Value = -96800  # N94: -46550 -> -96800

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 146286.28  # H40: 93900.37 -> 146286.28
$ws.Cells.Item(40, 9).Value = 170000.67  # I40: 113989.336 -> 170000.67
$ws.Cells.Item(40, 10).Value = 4000  # J40: 3500 -> 4000
$ws.Cells.Item(40, 11).Value = 170000.67  # K40: 113989.336 -> 170000.67
$ws.Cells.Item(40, 12).Value = 4000  # L40: 3500 -> 4000
$ws.Cells.Item(40, 13).Value = -169864.67  # M40: -113853.336 -> -169864.67
$ws.Cells.Item(40, 14).Value = -4272  # N40: -3772 -> -4272

$ws.Cells.Item(64, 8).Value = 30000  # H64: 13998 -> 30000
$ws.Cells.Item(64, 10).Value = 30000  # J64: 13998 -> 30000
$ws.Cells.Item(64, 12).Value = 30000  # L64: 13998 -> 30000
$ws.Cells.Item(64, 14).Value = -30450  # N64: -14448 -> -30450

$ws.Cells.Item(67, 8).Value = 30000  # H67: 13998 -> 30000
$ws.Cells.Item(67, 10).Value = 30000  # J67: 13998 -> 30000
$ws.Cells.Item(67, 12).Value = 30000  # L67: 13998 -> 30000
$ws.Cells.Item(67, 14).Value = -31560  # N67: -15558 -> -31560

$ws.Cells.Item(68, 8).Value = 1900  # H68: 0 -> 1900
$ws.Cells.Item(68, 9).Value = 1514.7059  # I68: 0 -> 1514.7059
$ws.Cells.Item(68, 10).Value = 8450  # J68: 0 -> 8450
$ws.Cells.Item(68, 11).Value = 1514.7059  # K68: 0 -> 1514.7059
$ws.Cells.Item(68, 12).Value = 8450  # L68: 0 -> 8450
$ws.Cells.Item(68, 13).Value = -765.7058999999999  # M68: None -> -765.7058999999999
$ws.Cells.Item(68, 14).Value = -9948  # N68: None -> -9948

$ws.Cells.Item(71, 8).Value = 1900  # H71: 0 -> 1900
$ws.Cells.Item(71, 9).Value = 1514.7059  # I71: 0 -> 1514.7059
$ws.Cells.Item(71, 10).Value = 8450  # J71: 0 -> 8450
$ws.Cells.Item(71, 11).Value = 7573.5295  # K71: 0 -> 7573.5295
$ws.Cells.Item(71, 12).Value = 42250  # L71: 0 -> 42250
$ws.Cells.Item(71, 13).Value = -3829.5295  # M71: None -> -3829.5295
$ws.Cells.Item(71, 14).Value = -49738  # N71: None -> -49738

$ws.Cells.Item(122, 8).Value = 3694.1177  # H122: 3753.125 -> 3694.1177
$ws.Cells.Item(122, 9).Value = 3078.5715  # I122: 3275 -> 3078.5715
$ws.Cells.Item(122, 10).Value = 4125  # J122: 3912.5 -> 4125
$ws.Cells.Item(122, 11).Value = 9235.7145  # K122: 9825 -> 9235.7145
$ws.Cells.Item(122, 12).Value = 12375  # L122: 11737.5 -> 12375
$ws.Cells.Item(122, 13).Value = -6785.7145  # M122: -7375 -> -6785.7145
$ws.Cells.Item(122, 14).Value = -17275  # N122: -16637.5 -> -17275

$ws.Cells.Item(132, 8).Value = 2798.3215  # H132: 3058.12 -> 2798.3215
$ws.Cells.Item(132, 9).Value = 2348.125  # I132: 2593.0952 -> 2348.125
$ws.Cells.Item(132, 11).Value = 7044.375  # K132: 7779.285600000001 -> 7044.375
$ws.Cells.Item(132, 13).Value = -4514.375  # M132: -5249.285600000001 -> -4514.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(27, 8).Value = 0  # H27: 40000 -> 0
$ws.Cells.Item(27, 10).Value = 0  # J27: 40000 -> 0
$ws.Cells.Item(27, 12).Value = 0  # L27: 40000 -> 0
$ws.Cells.Item(27, 14).ClearContents()  # N27: remove (was -40138)

$ws.Cells.Item(54, 8).Value = 9797  # H54: 19999 -> 9797
$ws.Cells.Item(54, 10).Value = 9797  # J54: 19999 -> 9797
$ws.Cells.Item(54, 12).Value = 9797  # L54: 19999 -> 9797
$ws.Cells.Item(54, 14).Value = -10837  # N54: -21039 -> -10837

$ws.Cells.Item(62, 8).Value = 5333.3335  # H62: 5500 -> 5333.3335
$ws.Cells.Item(62, 10).Value = 10000  # J62: 6750 -> 10000
$ws.Cells.Item(62, 12).Value = 10000  # L62: 6750 -> 10000
$ws.Cells.Item(62, 14).Value = -11248  # N62: -7998 -> -11248

$ws.Cells.Item(65, 8).Value = 5333.3335  # H65: 5500 -> 5333.3335
$ws.Cells.Item(65, 10).Value = 10000  # J65: 6750 -> 10000
$ws.Cells.Item(65, 12).Value = 50000  # L65: 33750 -> 50000
$ws.Cells.Item(65, 14).Value = -56240  # N65: -39990 -> -56240

$ws.Cells.Item(80, 8).Value = 73433.664  # H80: 60150.5 -> 73433.664
$ws.Cells.Item(80, 9).Value = 0  # I80: 70000 -> 0
$ws.Cells.Item(80, 10).Value = 73433.664  # J80: 50301 -> 73433.664
$ws.Cells.Item(80, 11).Value = 0  # K80: 70000 -> 0
$ws.Cells.Item(80, 12).Value = 73433.664  # L80: 50301 -> 73433.664
$ws.Cells.Item(80, 13).ClearContents()  # M80: remove (was -69002)
$ws.Cells.Item(80, 14).Value = -75429.664  # N80: -52297 -> -75429.664

$ws.Cells.Item(83, 8).Value = 73433.664  # H83: 60150.5 -> 73433.664
$ws.Cells.Item(83, 9).Value = 0  # I83: 70000 -> 0
$ws.Cells.Item(83, 10).Value = 73433.664  # J83: 50301 -> 73433.664
$ws.Cells.Item(83, 11).Value = 0  # K83: 210000 -> 0
$ws.Cells.Item(83, 12).Value = 220300.992  # L83: 150903 -> 220300.992
$ws.Cells.Item(83, 13).ClearContents()  # M83: remove (was -205008)
$ws.Cells.Item(83, 14).Value = -230284.992  # N83: -160887 -> -230284.992

$ws.Cells.Item(115, 8).Value = 39800  # H115: 32425 -> 39800
$ws.Cells.Item(115, 10).Value = 39800  # J115: 32425 -> 39800
$ws.Cells.Item(115, 12).Value = 39800  # L115: 32425 -> 39800
$ws.Cells.Item(115, 14).Value = -42934  # N115: -35559 -> -42934

$ws.Cells.Item(136, 8).Value = 2283.5535  # H136: 2434.25 -> 2283.5535
$ws.Cells.Item(136, 9).Value = 1839.75  # I136: 1962.1708 -> 1839.75
$ws.Cells.Item(136, 10).Value = 3910.8333  # J136: 4193.8184 -> 3910.8333
$ws.Cells.Item(136, 11).Value = 5519.25  # K136: 5886.512400000001 -> 5519.25
$ws.Cells.Item(136, 12).Value = 11732.4999  # L136: 12581.4552 -> 11732.4999
$ws.Cells.Item(136, 13).Value = -2969.25  # M136: -3336.512400000001 -> -2969.25
$ws.Cells.Item(136, 14).Value = -16832.4999  # N136: -17681.4552 -> -16832.4999
